$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value of 45190 (2023-09-21)
# for every data row (rows 2-498). Update it to 45192 (2023-09-23).
for ($row = 2; $row -le 498; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45190) {
        $cell.Value2 = 45192
    }
}
